$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed counts in-place (rows 2-6)
$ws.Range("B2").Value = 16
$ws.Range("B4").Value = 11
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 10
$ws.Range("A6").Value = 3

# Remove the old trailing rows 7-9 (data now ends at row 6), shifting rows up
$ws.Rows("7:9").Delete()
